$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.442728996276855
$ws.Range("B1").Value = 2.595495700836182
$ws.Range("C1").Value = 1.992204308509827
$ws.Range("D1").Value = 1.879971742630005
$ws.Range("E1").Value = 1.708329081535339
